$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("overall_results")
$ws2 = $wb.Worksheets.Item("bootstrap_results")

# --- overall_results: row 4 (mean) ---
$ws1.Range("C4").Value = 0.793
$ws1.Range("D4").Value = 0.781
$ws1.Range("E4").Value = 0.833
$ws1.Range("F4").Value = 0.882
$ws1.Range("G4").Value = 0.733
$ws1.Range("H4").Value = 0.717
$ws1.Range("I4").Value = 0.788
$ws1.Range("J4").Value = 0.8159999999999999
$ws1.Range("K4").Value = 0.747
$ws1.Range("L4").Value = 0.6909999999999999
$ws1.Range("M4").Value = 0.852
$ws1.Range("N4").Value = 0.82

# --- overall_results: row 5 (median) ---
$ws1.Range("C5").Value = 0.802
$ws1.Range("D5").Value = 0.756
$ws1.Range("E5").Value = 0.869
$ws1.Range("F5").Value = 0.887
$ws1.Range("G5").Value = 0.736
$ws1.Range("H5").Value = 0.719
$ws1.Range("I5").Value = 0.852
$ws1.Range("J5").Value = 0.8129999999999999
$ws1.Range("K5").Value = 0.758
$ws1.Range("L5").Value = 0.649
$ws1.Range("M5").Value = 0.928
$ws1.Range("N5").Value = 0.829

# --- overall_results: row 6 (std) ---
$ws1.Range("C6").Value = 0.049
$ws1.Range("D6").Value = 0.082
$ws1.Range("E6").Value = 0.131
$ws1.Range("F6").Value = 0.035
$ws1.Range("G6").Value = 0.067
$ws1.Range("H6").Value = 0.08799999999999999
$ws1.Range("I6").Value = 0.162
$ws1.Range("J6").Value = 0.057
$ws1.Range("K6").Value = 0.034
$ws1.Range("L6").Value = 0.083
$ws1.Range("M6").Value = 0.142
$ws1.Range("N6").Value = 0.031

# --- bootstrap_results: CI low / CI high columns (D, E) ---
$ws2.Range("E2").Value = 0.8057

$ws2.Range("D3").Value = 0.7578
$ws2.Range("E3").Value = 0.8026

$ws2.Range("D4").Value = 0.795
$ws2.Range("E4").Value = 0.8698

$ws2.Range("D5").Value = 0.8715000000000001
$ws2.Range("E5").Value = 0.8921

$ws2.Range("D6").Value = 0.737
$ws2.Range("E6").Value = 0.7564

$ws2.Range("D7").Value = 0.6699000000000001
$ws2.Range("E7").Value = 0.715

$ws2.Range("D8").Value = 0.8099
$ws2.Range("E8").Value = 0.8892

$ws2.Range("D9").Value = 0.8107
$ws2.Range("E9").Value = 0.8283

$ws2.Range("D10").Value = 0.7147
$ws2.Range("E10").Value = 0.7519

$ws2.Range("D11").Value = 0.694
$ws2.Range("E11").Value = 0.739

$ws2.Range("D12").Value = 0.7413
$ws2.Range("E12").Value = 0.8297

$ws2.Range("D13").Value = 0.8006
$ws2.Range("E13").Value = 0.8322000000000001
